# Auto-generated edit script for cryptos.xlsx data refresh
# Updates Price (D) and, for re-ranked coins, Coin (B) / Link (C) / Volume(1h) (E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '245.87'

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '22.02'

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.358'

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.05973'

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '3.395'

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '6.386'

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.8103'

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.9574'

# Row 10
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1431'
$ws.Range('E10').Value = '9WazirXWRX'

# Row 11
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07389'
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'

# Row 12
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.03377'
$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'

# Row 13
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.03066'
$ws.Range('E13').Value = '12BitrueCoinBTR'

# Row 14
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.09415'
$ws.Range('E14').Value = '13BitMartTokenBMX'

# Row 15
$ws.Range('B15').Value = 'MCDex'
$ws.Range('C15').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.996'
$ws.Range('E15').Value = '14MCDexMCB'

# Row 16
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.001592'
$ws.Range('E16').Value = '15BitForexTokenBF'

# Row 17
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.04812'
$ws.Range('E17').Value = '16CoinExTokenCET'

# Row 18
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0005913'
$ws.Range('E18').Value = '17OneONE'

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.006118'

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.005106'

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.0009880'

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.00006996'

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1333'

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0002463'

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.04145'

# Row 41
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.1073'
$ws.Range('E41').Value = '40BKEXTokenBKK'

# Row 42
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.002719'
$ws.Range('E42').Value = '41CEJICEJI'

# Row 43
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.003038'
$ws.Range('E43').Value = '42KickTokenKICK'

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.005818'

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00005255'

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.8505'
$ws.Range('E47').Value = '46CoinbaseStockTokenCOINBestin24h'

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.03231'

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.01011'
